$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text in A1 with the new daily conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldText = $ws1.Range("A1").Value()
$newText = $oldText -replace [regex]::Escape("1000 Bs = 9.76 = 40965.07 pesos"), "1000 Bs = 9.88 = 41578.03 pesos"
$newText = $newText -replace [regex]::Escape("40965.07 pesos = 9.68 = 942.83 Bs"), "41578.03 pesos = 9.83 = 959.34 Bs"
$ws1.Range("A1").Value = $newText

# --- Update the "tasas" sheet numeric cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 101.177
$ws2.Range("O10").Value = 4206.74

$ws2.Range("N12").Value = 4229.99
$ws2.Range("O12").Value = 97.6
